$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) sometimes holds values that look like plain numbers
# (e.g. "590.72"). Excel's Value setter auto-coerces those into doubles,
# which both loses the original text formatting and introduces binary
# floating point noise. Force the whole Price column to Text first so every
# assignment below is stored verbatim as a string, then restore the
# worksheet's default style so the cells end up styled exactly as they
# started (no stray per-cell "s" attribute left behind).
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value  = "67.045.62"
$ws.Range("E2").Value  = "  -0.47%  "

$ws.Range("D3").Value  = "2.612.08"
$ws.Range("E3").Value  = "  -1.19%  "

$ws.Range("E4").Value  = "  -0.03%  "

$ws.Range("D5").Value  = "590.72"
$ws.Range("E5").Value  = "  -1.34%  "

$ws.Range("D6").Value  = "165.59"
$ws.Range("E6").Value  = "  -0.41%  "

$ws.Range("E7").Value  = "  +0.00%  "

$ws.Range("E8").Value  = "  -2.19%  "

$ws.Range("D9").Value  = "2.608.83"
$ws.Range("E9").Value  = "  -1.31%  "

$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  -4.71%  "

$ws.Range("D12").Value = "0.363"

$ws.Range("E13").Value = "  -0.69%  "

$ws.Range("D14").Value = "27.32"
$ws.Range("E14").Value = "  -2.49%  "

$ws.Range("D15").Value = "3.084.19"

$ws.Range("E16").Value = "  -2.50%  "

$ws.Range("D17").Value = "66.905.75"
$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("D18").Value = "2.624.00"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").Value = "11.79"
$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("D20").Value = "7.82"
$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D21").Value = "354.29"
$ws.Range("E21").Value = "  -2.73%  "

$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  -2.84%  "

$ws.Range("D23").Value = "4.63"
$ws.Range("E23").Value = "  -3.38%  "

$ws.Range("D24").Value = "10.56"
$ws.Range("E24").Value = "  -5.03%  "

$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("E26").Value = "  -4.70%  "

$ws.Range("D27").Value = "68.99"
$ws.Range("E27").Value = "  -2.84%  "

$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").Value = "0.0₃0996"
$ws.Range("E30").Value = "  -2.72%  "

$ws.Range("D31").Value = "543.77"
$ws.Range("E31").Value = "  -2.04%  "

$ws.Range("D32").Value = "7.87"
$ws.Range("E32").Value = "  -2.12%  "

$ws.Range("E33").Value = "  -3.76%  "

$ws.Range("E34").Value = "  -2.78%  "

$ws.Range("E35").Value = "  +0.79%  "

$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("E37").Value = "  -3.67%  "

$ws.Range("D38").Value = "157.23"
$ws.Range("E38").Value = "  -0.28%  "

$ws.Range("D39").Value = "18.93"
$ws.Range("E39").Value = "  -2.56%  "

$ws.Range("E40").Value = "  -2.12%  "

$ws.Range("E41").Value = "  +1.78%  "

$ws.Range("D42").Value = "1.81"
$ws.Range("E42").Value = "  -1.12%  "

$ws.Range("D43").Value = "5.14"
$ws.Range("E43").Value = "  -2.48%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  -4.64%  "

$ws.Range("D46").Value = "0.0₆0299"
$ws.Range("E46").Value = "  -1.54%  "

$ws.Range("D47").Value = "151.51"
$ws.Range("E47").Value = "  -1.74%  "

$ws.Range("D48").Value = "0.575"
$ws.Range("E48").Value = "  -3.66%  "

$ws.Range("D49").Value = "3.77"
$ws.Range("E49").Value = "  -3.04%  "

$ws.Range("E50").Value = "  -2.10%  "

$ws.Range("E51").Value = "  -1.18%  "

# Restore the Price column's style to the sheet default now that every
# value has been committed as text.
$priceCol.Style = "Normal"
